$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, Fecha(serial), Volumen, PrecioMin, PrecioMax, PrecioPromPonderado, Precio$/Kg
$data = @(
    @(2, 44377, 40, 14000, 15000, 14500, 1115),
    @(3, 44362, 40, 15000, 16000, 15500, 1192),
    @(4, 44488, 40, 16000, 17000, 16500, 1269),
    @(5, 44691, 100, 12000, 13000, 12500, 962),
    @(6, 44453, 50, 14000, 15000, 14600, 1123),
    @(7, 44350, 40, 23000, 25000, 24000, 1846),
    @(8, 44435, 100, 13000, 14000, 13500, 1038),
    @(9, 44610, 50, 17000, 18000, 17400, 1338),
    @(10, 44664, 50, 11000, 12000, 11600, 892),
    @(11, 44425, 60, 14000, 15000, 14500, 1115),
    @(12, 44355, 60, 18000, 20000, 19000, 1462),
    @(13, 44159, 60, 30000, 32000, 31000, 2385),
    @(14, 44503, 35, 15000, 16000, 15429, 1187),
    @(15, 44509, 100, 15000, 16000, 15500, 1192),
    @(16, 44523, 40, 15000, 16000, 15500, 1192),
    @(17, 44334, 50, 26000, 28000, 27200, 2092),
    @(18, 44313, 50, 25000, 26000, 25600, 1969),
    @(19, 44467, 100, 13000, 14000, 13500, 1038),
    @(20, 44462, 60, 14000, 15000, 14500, 1115),
    @(21, 44510, 40, 15000, 16000, 15500, 1192),
    @(22, 44320, 50, 26000, 28000, 26800, 2062),
    @(23, 44264, 40, 30000, 32000, 31000, 2385),
    @(24, 44316, 50, 27000, 28000, 27400, 2108),
    @(25, 44433, 100, 13000, 14000, 13500, 1038),
    @(26, 44308, 50, 26000, 27000, 26400, 2031),
    @(27, 44383, 50, 15000, 16000, 15400, 1185),
    @(28, 44327, 50, 24000, 25000, 24400, 1877),
    @(29, 44474, 40, 13000, 14000, 13500, 1038)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value2 = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}
